$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(56).Insert()

$ws.Cells.Item(56, 1).Value = 9
$ws.Cells.Item(56, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(56, 3).Value = "Metropolitana"
$ws.Cells.Item(56, 4).Value = (Get-Date -Year 2021 -Month 11 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(56, 5).Value = 13
$ws.Cells.Item(56, 6).Value = 300000001
$ws.Cells.Item(56, 7).Value = "Rabanito"
$ws.Cells.Item(56, 8).Value = "Sin especificar"
$ws.Cells.Item(56, 9).Value = "Primera"
$ws.Cells.Item(56, 10).Value = 7900
$ws.Cells.Item(56, 11).Value = 2500
$ws.Cells.Item(56, 12).Value = 3000
$ws.Cells.Item(56, 13).Value = 2747
$ws.Cells.Item(56, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(56, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(56, 16).Value = 27
$ws.Cells.Item(56, 17).Value = 100
$ws.Cells.Item(56, 18).Value = "Hortaliza"
